$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value while keeping the cell a genuine
# text cell (matches the workbook's existing "inline string" price/hour
# columns) and without leaving any lingering custom number-format style
# behind on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "249.71"
Set-TextValue "G2" "15"

# Row 3
Set-TextValue "D3" "24.11"
Set-TextValue "G3" "15"

# Row 4
Set-TextValue "D4" "5.956"
Set-TextValue "G4" "15"

# Row 5
Set-TextValue "G5" "15"

# Row 6
Set-TextValue "D6" "3.428"
Set-TextValue "G6" "15"

# Row 7
Set-TextValue "D7" "6.530"
Set-TextValue "G7" "15"

# Row 8
Set-TextValue "G8" "15"

# Row 9
Set-TextValue "D9" "0.7978"
Set-TextValue "G9" "15"

# Row 10
Set-TextValue "D10" "0.1485"
Set-TextValue "G10" "15"

# Row 11
Set-TextValue "D11" "0.07806"
Set-TextValue "G11" "15"

# Row 12
Set-TextValue "D12" "0.03302"
Set-TextValue "G12" "15"

# Row 13
Set-TextValue "D13" "0.03000"
Set-TextValue "G13" "15"

# Row 14
Set-TextValue "D14" "0.09241"
Set-TextValue "G14" "15"

# Row 15
Set-TextValue "D15" "3.563"
Set-TextValue "G15" "15"

# Row 16
Set-TextValue "D16" "0.001666"
Set-TextValue "G16" "15"

# Row 17
Set-TextValue "D17" "0.04759"
Set-TextValue "G17" "15"

# Row 18 (was TigerCash -> now One)
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0006030"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue "G18" "15"

# Row 19 (was HotbitToken -> now TigerCash)
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D19" "0.006212"
$ws.Range("E19").Value = "18TigerCashTCH"
Set-TextValue "G19" "15"

# Row 20 (was BitKan -> now HotbitToken)
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D20" "0.005573"
$ws.Range("E20").Value = "19HotbitTokenHTB"
Set-TextValue "G20" "15"

# Row 21 (was NitroEx -> now BitKan)
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D21" "0.001065"
$ws.Range("E21").Value = "20BitKanKAN"
Set-TextValue "G21" "15"

# Row 22 (was LEO -> now NitroEx)
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D22" "0.0001501"
$ws.Range("E22").Value = "21NitroExNTX"
Set-TextValue "G22" "15"

# Row 23 (was BTSEToken -> now LEO)
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D23" "3.699"
$ws.Range("E23").Value = "22LEOLEO"
Set-TextValue "G23" "15"

# Row 24 (was One -> now BTSEToken)
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D24" "2.212"
$ws.Range("E24").Value = "23BTSETokenBTSE"
Set-TextValue "G24" "15"

# Row 25
Set-TextValue "D25" "0.3358"
Set-TextValue "G25" "15"

# Row 26
Set-TextValue "D26" "0.1253"
Set-TextValue "G26" "15"

# Row 27
Set-TextValue "D27" "0.0006475"
Set-TextValue "G27" "15"

# Row 28
Set-TextValue "G28" "15"

# Row 29
Set-TextValue "G29" "15"

# Row 30
Set-TextValue "G30" "15"

# Row 31
Set-TextValue "G31" "15"

# Row 32
Set-TextValue "G32" "15"

# Row 33
Set-TextValue "G33" "15"

# Row 34
Set-TextValue "G34" "15"

# Row 35
Set-TextValue "G35" "15"

# Row 36
Set-TextValue "G36" "15"

# Row 37
Set-TextValue "G37" "15"

# Row 38
Set-TextValue "G38" "15"

# Row 39
Set-TextValue "G39" "15"

# Row 40
Set-TextValue "D40" "0.04405"
Set-TextValue "G40" "15"

# Row 41
Set-TextValue "D41" "0.007034"
Set-TextValue "G41" "15"

# Row 42 (was BKEXToken -> now CEJI)
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003601"
$ws.Range("E42").Value = "41CEJICEJI"
Set-TextValue "G42" "15"

# Row 43 (was CEJI -> now BKEXToken)
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1066"
$ws.Range("E43").Value = "42BKEXTokenBKK"
Set-TextValue "G43" "15"

# Row 44
Set-TextValue "D44" "0.009773"
Set-TextValue "G44" "15"

# Row 45
Set-TextValue "D45" "0.002460"
$ws.Range("E45").Value = "44ACDXExchangeACXTBestin24h"
Set-TextValue "G45" "15"

# Row 46
Set-TextValue "D46" "0.00005894"
Set-TextValue "G46" "15"

# Row 47
Set-TextValue "D47" "0.00000000750"
Set-TextValue "G47" "15"

# Row 48
Set-TextValue "D48" "0.9901"
Set-TextValue "G48" "15"

# Row 49
$ws.Range("E49").Value = "48BOLOBOLO"
Set-TextValue "G49" "15"

# Row 50
Set-TextValue "D50" "0.00002101"
Set-TextValue "G50" "15"

# Row 51
Set-TextValue "D51" "0.01010"
Set-TextValue "G51" "15"
